{"js": "const body = context.document.body;\nconst results = body.search(\"Version 2.\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"Version 1.\", \"Replace\");\n} else {\n  throw new Error(\"Text 'Version 2.' not found in document body.\");\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# \"Versi\" + \"on\" were historically split across two runs with identical\n# formatting; re-running a no-op Find/Replace over the whole word merges\n# them into a single run, matching the canonical \"Version\" run.\n$find1 = $d.Content.Find\n$find1.Execute(\n  \"Version\",    # FindText\n  $true,        # MatchCase\n  $false,       # MatchWholeWord\n  $false,       # MatchWildcards\n  $false,       # MatchSoundsLike\n  $false,       # MatchAllWordForms\n  $true,        # Forward\n  1,            # Wrap (wdFindContinue)\n  $false,       # Format\n  \"Version\",    # ReplaceWith\n  2             # Replace (wdReplaceAll)\n) | Out-Null\n\n# Narrow replace of just the digit (keeps the _GoBack bookmark and the\n# spellEnd proofErr mark untouched, since the matched span doesn't cross\n# them).\n$find2 = $d.Content.Find\n$find2.Execute(\n  \" 2\",         # FindText\n  $true, $false, $false, $false, $false,\n  $true, 1, $false,\n  \" 1\",         # ReplaceWith\n  2\n) | Out-Null\n\n# The trailing \".\" used to sit in its own run AFTER the _GoBack bookmark.\n# Drop it, then insert a fresh \".\" immediately before the bookmark so it\n# becomes part of the \" 1\" run (\" 1.\") and the bookmark ends up right\n# after it again - matching \"Version 1.\" with the bookmark at the tail.\n$full = $d.Content\n$periodRange = $d.Range($full.End - 2, $full.End - 1)\nif ($periodRange.Text -eq \".\") {\n  $periodRange.Text = \"\"\n}\n\n$bm = $d.Bookmarks(\"_GoBack\")\n$bm.Range.InsertBefore(\".\")\n"}
